$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Delete the empty paragraph that sits right before the "BootStrat"
#    Heading2 paragraph.
# ------------------------------------------------------------------
$headingRng = $d.Content.Duplicate
$headingRng.Find.Execute("BootStrat", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$headingStart = $headingRng.Start

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.End -eq $headingStart -and $p.Range.Text.Length -le 1) {
        $p.Range.Delete()
        break
    }
}

# ------------------------------------------------------------------
# 2) Remove the old "_GoBack" bookmark that used to sit around
#    "Include via CDN". Do this *before* adding the new one below,
#    since only one bookmark with a given name can live in the
#    document's bookmark index at a time.
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 3) Turn "BootStrat" into two runs "BootStra" + "p" (i.e. fix the
#    typo to "BootStrap") and add a fresh "_GoBack" bookmark right
#    before those runs, inside the Heading2 paragraph.
# ------------------------------------------------------------------
$rng2 = $d.Content.Duplicate
$rng2.Find.Execute("BootStrat", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$start2 = $rng2.Start
$end2 = $rng2.End
$target2 = $d.Range($start2, $end2)
$xml2 = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t>BootStra</w:t></w:r><w:r><w:t>p</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$target2.InsertXML($xml2)

# ------------------------------------------------------------------
# 4) Merge the hyperlink's three runs ("https://getboot" + "s" +
#    "trap.com/") into a single run "https://getbootstrap.com/".
#    A same-text replace is a no-op for this engine, so first replace
#    with a distinct placeholder of identical length (forcing a true
#    merge of the backing runs), then replace that placeholder with
#    the real text.
# ------------------------------------------------------------------
$rng3 = $d.Content.Duplicate
$rng3.Find.Execute("https://getbootstrap.com/", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$start3 = $rng3.Start
$end3 = $rng3.End
$len3 = $end3 - $start3
$placeholder = ""
for ($j = 0; $j -lt $len3; $j++) { $placeholder += "Z" }
$target3 = $d.Range($start3, $end3)
$target3.Text = $placeholder
$target3b = $d.Range($start3, $start3 + $len3)
$target3b.Text = "https://getbootstrap.com/"
